$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2023-2024")

$ws.Range("B8").Value = "MPAL"
$ws.Range("C8").Value = "TP"
$ws.Range("D8").Value = "X"
$ws.Range("G8").Value = "45min de prise en main de Issues, Classroom et Projects, 30min de rédaction des US."
$ws.Range("I8").Value = "Difficulté à différencier le point d'entrée Classroom vs le dépôt Github."

$ws.Range("I7").Select()
